$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value2 = 1200
$ws.Range("I16").Value2 = 0
$ws.Range("J16").Value2 = 1200
$ws.Range("K16").Value2 = 0
$ws.Range("L16").Value2 = 1200
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value2 = -1660
$ws.Range("H17").Value2 = 908333.0600000001
$ws.Range("I17").Value2 = 436.32144
$ws.Range("J17").Value2 = 1397200.5
$ws.Range("K17").Value2 = 1308.96432
$ws.Range("L17").Value2 = 4191601.5
$ws.Range("M17").Value2 = -1140.96432
$ws.Range("N17").Value2 = -4191937.5
$ws.Range("H100").Value2 = 1905.2084
$ws.Range("I100").Value2 = 1072.0588
$ws.Range("J100").Value2 = 3928.5715
$ws.Range("K100").Value2 = 1072.0588
$ws.Range("L100").Value2 = 3928.5715
$ws.Range("M100").Value2 = -531.0588
$ws.Range("N100").Value2 = -5010.5715
$ws.Range("H132").Value2 = 2427.627
$ws.Range("I132").Value2 = 2339.0168
$ws.Range("J132").Value2 = 3081.125
$ws.Range("K132").Value2 = 7017.0504
$ws.Range("L132").Value2 = 9243.375
$ws.Range("M132").Value2 = -4487.0504
$ws.Range("N132").Value2 = -14303.375
$ws.Range("H133").Value2 = 89780
$ws.Range("J133").Value2 = 89780
$ws.Range("L133").Value2 = 89780
$ws.Range("N133").Value2 = -99900
$ws.Range("H137").Value2 = 1255.3662
$ws.Range("I137").Value2 = 1143.122
$ws.Range("J137").Value2 = 1408.7667
$ws.Range("K137").Value2 = 3429.366
$ws.Range("L137").Value2 = 4226.300099999999
$ws.Range("M137").Value2 = -879.366
$ws.Range("N137").Value2 = -9326.3001
$ws.Range("H138").Value2 = 1651.6
$ws.Range("I138").Value2 = 841.8570999999999
$ws.Range("J138").Value2 = 1966.5
$ws.Range("K138").Value2 = 2525.5713
$ws.Range("L138").Value2 = 5899.5
$ws.Range("M138").Value2 = 2614.4287
$ws.Range("N138").Value2 = -16179.5
$ws.Range("H141").Value2 = 2392.8223
$ws.Range("I141").Value2 = 965.6585
$ws.Range("J141").Value2 = 17021.25
$ws.Range("K141").Value2 = 2896.9755
$ws.Range("L141").Value2 = 51063.75
$ws.Range("M141").Value2 = 2283.0245
$ws.Range("N141").Value2 = -61423.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value2 = 5940.57
$ws.Range("I32").Value2 = 5454.25
$ws.Range("J32").Value2 = 9506.916999999999
$ws.Range("K32").Value2 = 5454.25
$ws.Range("L32").Value2 = 9506.916999999999
$ws.Range("M32").Value2 = -5167.25
$ws.Range("N32").Value2 = -10080.917
$ws.Range("H61").Value2 = 4387420
$ws.Range("I61").Value2 = 5377583
$ws.Range("J61").Value2 = 2412.2856
$ws.Range("K61").Value2 = 5377583
$ws.Range("L61").Value2 = 2412.2856
$ws.Range("M61").Value2 = -5377371
$ws.Range("N61").Value2 = -2836.2856
$ws.Range("H74").Value2 = 1844.8448
$ws.Range("I74").Value2 = 865.9643
$ws.Range("J74").Value2 = 2758.4666
$ws.Range("K74").Value2 = 865.9643
$ws.Range("L74").Value2 = 2758.4666
$ws.Range("M74").Value2 = 8.03570000000002
$ws.Range("N74").Value2 = -4506.4666
$ws.Range("H77").Value2 = 1844.8448
$ws.Range("I77").Value2 = 865.9643
$ws.Range("J77").Value2 = 2758.4666
$ws.Range("K77").Value2 = 4329.8215
$ws.Range("L77").Value2 = 13792.333
$ws.Range("M77").Value2 = 38.17849999999999
$ws.Range("N77").Value2 = -22528.333
$ws.Range("H102").Value2 = 1505
$ws.Range("I102").Value2 = 1505
$ws.Range("J102").Value2 = 0
$ws.Range("K102").Value2 = 1505
$ws.Range("L102").Value2 = 0
$ws.Range("M102").Value2 = 117
$ws.Range("N102").ClearContents()
$ws.Range("H110").Value2 = 48327.79
$ws.Range("I110").Value2 = 60978.535
$ws.Range("J110").Value2 = 887.5
$ws.Range("K110").Value2 = 60978.535
$ws.Range("L110").Value2 = 887.5
$ws.Range("M110").Value2 = -58933.535
$ws.Range("N110").Value2 = -4977.5
$ws.Range("H132").Value2 = 3703.93
$ws.Range("I132").Value2 = 3501.0833
$ws.Range("J132").Value2 = 4051.6667
$ws.Range("K132").Value2 = 10503.2499
$ws.Range("L132").Value2 = 12155.0001
$ws.Range("M132").Value2 = -7973.249899999999
$ws.Range("N132").Value2 = -17215.0001
$ws.Range("H136").Value2 = 4387420
$ws.Range("I136").Value2 = 5377583
$ws.Range("J136").Value2 = 2412.2856
$ws.Range("K136").Value2 = 16132749
$ws.Range("L136").Value2 = 7236.8568
$ws.Range("M136").Value2 = -16130199
$ws.Range("N136").Value2 = -12336.8568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value2 = 529.61536
$ws.Range("I64").Value2 = 519.5
$ws.Range("J64").Value2 = 545.8
$ws.Range("K64").Value2 = 519.5
$ws.Range("L64").Value2 = 545.8
$ws.Range("M64").Value2 = -294.5
$ws.Range("N64").Value2 = -995.8
$ws.Range("H67").Value2 = 529.61536
$ws.Range("I67").Value2 = 519.5
$ws.Range("J67").Value2 = 545.8
$ws.Range("K67").Value2 = 519.5
$ws.Range("L67").Value2 = 545.8
$ws.Range("M67").Value2 = 260.5
$ws.Range("N67").Value2 = -2105.8
$ws.Range("H134").Value2 = 1649.0878
$ws.Range("I134").Value2 = 1550.3062
$ws.Range("K134").Value2 = 4650.9186
$ws.Range("M134").Value2 = -2115.9186

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value2 = 24143.715
$ws.Range("J4").Value2 = 24143.715
$ws.Range("L4").Value2 = 24143.715
$ws.Range("N4").Value2 = -24367.715
$ws.Range("H31").Value2 = 4070.5144
$ws.Range("I31").Value2 = 1315.3235
$ws.Range("J31").Value2 = 6672.6387
$ws.Range("K31").Value2 = 1315.3235
$ws.Range("L31").Value2 = 6672.6387
$ws.Range("M31").Value2 = -1020.3235
$ws.Range("N31").Value2 = -7262.6387
$ws.Range("H34").Value2 = 4070.5144
$ws.Range("I34").Value2 = 1315.3235
$ws.Range("J34").Value2 = 6672.6387
$ws.Range("K34").Value2 = 1315.3235
$ws.Range("L34").Value2 = 6672.6387
$ws.Range("M34").Value2 = -1113.3235
$ws.Range("N34").Value2 = -7076.6387
$ws.Range("H58").Value2 = 1115.4814
$ws.Range("I58").Value2 = 812.1667
$ws.Range("J58").Value2 = 1722.1111
$ws.Range("K58").Value2 = 812.1667
$ws.Range("L58").Value2 = 1722.1111
$ws.Range("M58").Value2 = -609.1667
$ws.Range("N58").Value2 = -2128.1111
$ws.Range("H134").Value2 = 2799.1128
$ws.Range("I134").Value2 = 2881.3542
$ws.Range("J134").Value2 = 2517.1428
$ws.Range("K134").Value2 = 8644.062600000001
$ws.Range("L134").Value2 = 7551.428400000001
$ws.Range("M134").Value2 = -6109.062600000001
$ws.Range("N134").Value2 = -12621.4284
$ws.Range("H136").Value2 = 1115.4814
$ws.Range("I136").Value2 = 812.1667
$ws.Range("J136").Value2 = 1722.1111
$ws.Range("K136").Value2 = 2436.5001
$ws.Range("L136").Value2 = 5166.3333
$ws.Range("M136").Value2 = 113.4998999999998
$ws.Range("N136").Value2 = -10266.3333
$ws.Range("H141").Value2 = 38694.26
$ws.Range("J141").Value2 = 31362.182
$ws.Range("L141").Value2 = 31362.182
$ws.Range("N141").Value2 = -41722.182

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value2 = 1396.3143
$ws.Range("J5").Value2 = 2074.7
$ws.Range("L5").Value2 = 6224.099999999999
$ws.Range("N5").Value2 = -6448.099999999999
$ws.Range("H113").Value2 = 480.81357
$ws.Range("I113").Value2 = 479.77777
$ws.Range("J113").Value2 = 481.6875
$ws.Range("K113").Value2 = 1439.33331
$ws.Range("L113").Value2 = 1445.0625
$ws.Range("M113").Value2 = 730.66669
$ws.Range("N113").Value2 = -5785.0625
$ws.Range("H128").Value2 = 200000
$ws.Range("I128").Value2 = 200000
$ws.Range("K128").Value2 = 600000
$ws.Range("M128").Value2 = -595020
$ws.Range("H131").Value2 = 3684
$ws.Range("I131").Value2 = 497.66666
$ws.Range("J131").Value2 = 4503.343
$ws.Range("K131").Value2 = 1492.99998
$ws.Range("L131").Value2 = 13510.029
$ws.Range("M131").Value2 = 3547.00002
$ws.Range("N131").Value2 = -23590.029
$ws.Range("H135").Value2 = 1396.3143
$ws.Range("J135").Value2 = 2074.7
$ws.Range("L135").Value2 = 18672.3
$ws.Range("N135").Value2 = -23742.3

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H55").Value2 = 10250
$ws.Range("I55").Value2 = 500
$ws.Range("J55").Value2 = 20000
$ws.Range("K55").Value2 = 500
$ws.Range("L55").Value2 = 20000
$ws.Range("M55").Value2 = -173
$ws.Range("N55").Value2 = -20654
$ws.Range("H132").Value2 = 2271.8428
$ws.Range("I132").Value2 = 1845.3208
$ws.Range("J132").Value2 = 3601.5881
$ws.Range("K132").Value2 = 5535.9624
$ws.Range("L132").Value2 = 10804.7643
$ws.Range("M132").Value2 = -3005.9624
$ws.Range("N132").Value2 = -15864.7643
$ws.Range("H139").Value2 = 220145.2
$ws.Range("J139").Value2 = 220145.2
$ws.Range("L139").Value2 = 220145.2
$ws.Range("N139").Value2 = -230425.2
$ws.Range("H141").Value2 = 78000
$ws.Range("J141").Value2 = 78000
$ws.Range("L141").Value2 = 78000
$ws.Range("N141").Value2 = -88360

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value2 = 6000501.5
$ws.Range("I2").Value2 = 1001
$ws.Range("K2").Value2 = 1001
$ws.Range("M2").Value2 = -889
$ws.Range("H68").Value2 = 1453.8182
$ws.Range("I68").Value2 = 1419.5918
$ws.Range("J68").Value2 = 1733.3334
$ws.Range("K68").Value2 = 1419.5918
$ws.Range("L68").Value2 = 1733.3334
$ws.Range("M68").Value2 = -670.5917999999999
$ws.Range("N68").Value2 = -3231.3334
$ws.Range("H71").Value2 = 1453.8182
$ws.Range("I71").Value2 = 1419.5918
$ws.Range("J71").Value2 = 1733.3334
$ws.Range("K71").Value2 = 7097.959
$ws.Range("L71").Value2 = 8666.666999999999
$ws.Range("M71").Value2 = -3353.959
$ws.Range("N71").Value2 = -16154.667
$ws.Range("H132").Value2 = 4200.35
$ws.Range("I132").Value2 = 3881.818
$ws.Range("J132").Value2 = 4589.6665
$ws.Range("K132").Value2 = 11645.454
$ws.Range("L132").Value2 = 13768.9995
$ws.Range("M132").Value2 = -9115.454000000002
$ws.Range("N132").Value2 = -18828.9995
$ws.Range("H136").Value2 = 2977450
$ws.Range("I136").Value2 = 1052.8206
$ws.Range("J136").Value2 = 9805655
$ws.Range("K136").Value2 = 3158.4618
$ws.Range("L136").Value2 = 29416965
$ws.Range("M136").Value2 = -608.4618
$ws.Range("N136").Value2 = -29422065
$ws.Range("H140").Value2 = 40183.453
$ws.Range("J140").Value2 = 40183.453
$ws.Range("L140").Value2 = 40183.453
$ws.Range("N140").Value2 = -50543.453

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value2 = 160714
$ws.Range("J46").Value2 = 160714
$ws.Range("L46").Value2 = 160714
$ws.Range("N46").Value2 = -161176
$ws.Range("H132").Value2 = 6207528.5
$ws.Range("I132").Value2 = 2034.1724
$ws.Range("J132").Value2 = 16205270
$ws.Range("K132").Value2 = 6102.5172
$ws.Range("L132").Value2 = 48615810
$ws.Range("M132").Value2 = -3572.5172
$ws.Range("N132").Value2 = -48620870
$ws.Range("H134").Value2 = 160714
$ws.Range("J134").Value2 = 160714
$ws.Range("L134").Value2 = 482142
$ws.Range("N134").Value2 = -487212
$ws.Range("H136").Value2 = 1175.02
$ws.Range("I136").Value2 = 1130.5264
$ws.Range("J136").Value2 = 1315.9166
$ws.Range("K136").Value2 = 3391.5792
$ws.Range("L136").Value2 = 3947.7498
$ws.Range("M136").Value2 = -841.5792000000001
$ws.Range("N136").Value2 = -9047.7498
$ws.Range("H140").Value2 = 80000
$ws.Range("J140").Value2 = 80000
$ws.Range("L140").Value2 = 80000
$ws.Range("N140").Value2 = -90360
$ws.Range("H141").Value2 = 51375
$ws.Range("J141").Value2 = 51375
$ws.Range("L141").Value2 = 51375
$ws.Range("N141").Value2 = -61735

